## Design decisions tracking.xlsx - content update
##
## Adds the new "bidding / watchlist / my-bids" design-decision rows to the
## R4 and R5 sheets, and leaves the cursor/zoom/active-tab the way the
## author left the workbook (R4 selected, cell B11 active on R4, cell C6
## active on R5).
##
## Shared strings are written in the same order the author appears to have
## typed them in (interleaving R4 and R5), so the resulting shared-string
## table lines up with the authored file as closely as possible.

$wb = $excel.ActiveWorkbook

$sheetR4 = $wb.Worksheets.Item(4)   # "R4"
$sheetR5 = $wb.Worksheets.Item(5)   # "R5"

# --- R4 row 8 -----------------------------------------------------------
$sheetR4.Range("B8").Value = "The bid form is only available to buyers who are logged in, and only before the auction end date."
$sheetR4.Range("C8").Value = "This is to prevent sellers and users who are not logged in from making a bid."

# --- R4 row 7 -----------------------------------------------------------
$sheetR4.Range("B7").Value = "When a bid is placed, we check if the auction end date is in the past."
$sheetR4.Range("C7").Value = "This it to prevent users from loading the listing page, waiting for the auction to end, and then submitting a bid."

# --- R5 rows 5-7 ---------------------------------------------------------
$sheetR5.Range("B5").Value = "Number of watchers is shown in the listing page, and updated live when the user adds or removes the listing from their watchlist"
$sheetR5.Range("B6").Value = "When a bid is placed on an auction that the buyer doesn't have on their wathclist, they are prompted to add to watchlist"
$sheetR5.Range("B7").Value = "When a bid is placed, all buyers who have the auction on their watchlist receive an email notification. Only the previous highest bidder is notified that they were outbid."

# --- R4 row 5 -------------------------------------------------------------
$sheetR4.Range("B5").Value = "Summary information is shown for each auction. This varies depending on whether the auction has ended or not."

# --- R4 row 9 -------------------------------------------------------------
$sheetR4.Range("B9").Value = "Bid history (bid amount, bid date, and bid username) for each auction is visible to all users."

# --- R4 row 5 (col C) ------------------------------------------------------
$sheetR4.Range("C5").Value = "Auctions that have ended have more limited information - they exlude starting price and minimum increment, but include final price."

# --- R4 row 9 (col C) ------------------------------------------------------
$sheetR4.Range("C9").Value = "All users are able to see bid history, including those who are not logged in."

# --- R4 row 6 --------------------------------------------------------------
$sheetR4.Range("B6").Value = "When a bid is placed, we verify the bid amount and format - only numbers with up to 2 decimal points allowed, and only numbers >= minimum bid are allowed (min. bid determined by current price and min. increment)"
$sheetR4.Range("C6").Value = "The bid input is quite crucial to the whole website, so we validate it both on the front end through the bid form, and in the back end through a SQL function."

# --- R4 row 10 --------------------------------------------------------------
$sheetR4.Range("B10").Value = "Buyers have a My Bids tab, where they can see, filter and sort all auctions they've bid on, and all individual bids. They also see what the current highest bid is, as well as who's winning the auction."

# --- R5 row 8 -----------------------------------------------------------
$sheetR5.Range("B8").Value = "Buyers have a My Watchlist tab, where they can see all their watched auctions, filter and sort them, as well as remove auctions from watchlist. They also see the current highest bid and bidder."
$sheetR5.Range("C8").Value = "This page has a lot of overlap with the My Bids page, but the key difference is that buyers can use this to track auctions before deciding whether to bid or not."

# --- R4 row 11 --------------------------------------------------------------
$sheetR4.Range("B11").Value = "When an auction ends, both buyers and sellers can see the outcome through the My Bids/My Watchlist and My Listings pages respectively. In addition, sellers receive an email notifying them of their auction's outcome, and winning buyers will also receive an email notification."

# --- view state: leave R5 showing C6 selected at 100% zoom ------------------
$sheetR5.Activate()
$sheetR5.Range("C6").Select()
$excel.ActiveWindow.Zoom = 100

# --- view state: leave R4 as the active/selected tab, B11 selected, 100% zoom
$sheetR4.Activate()
$sheetR4.Range("B11").Select()
$excel.ActiveWindow.Zoom = 100
